# Apply updated scraped values (想去人数 / 最低票价) per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 27001
$ws1.Range("G3").Value = "已售罄"
$ws1.Range("F5").Value = 637
$ws1.Range("F6").Value = 185
$ws1.Range("F7").Value = 561
$ws1.Range("F9").Value = 373
$ws1.Range("F11").Value = 195
$ws1.Range("F12").Value = 54
$ws1.Range("F13").Value = 312
$ws1.Range("F14").Value = 94
$ws1.Range("F15").Value = 469
$ws1.Range("F17").Value = 1610
$ws1.Range("F18").Value = 244
$ws1.Range("F19").Value = 443
$ws1.Range("G19").Value = "不可售"
$ws1.Range("F20").Value = 133
$ws1.Range("F21").Value = 453
$ws1.Range("F22").Value = 7
$ws1.Range("F23").Value = 109
$ws1.Range("F24").Value = 121

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 40
$ws2.Range("F9").Value = 2
$ws2.Range("F11").Value = 447
$ws2.Range("F17").Value = 73
$ws2.Range("F24").Value = 20

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5165
$ws3.Range("F3").Value = 266

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 5165
$ws4.Range("F4").Value = 266
$ws4.Range("F5").Value = 27001
$ws4.Range("G5").Value = "已售罄"
$ws4.Range("F9").Value = 637
$ws4.Range("F12").Value = 185
$ws4.Range("F15").Value = 40
$ws4.Range("F16").Value = 2
$ws4.Range("F18").Value = 447
$ws4.Range("F19").Value = 561
$ws4.Range("F23").Value = 373
$ws4.Range("F25").Value = 195
$ws4.Range("F26").Value = 54
$ws4.Range("F28").Value = 312
$ws4.Range("F29").Value = 94
$ws4.Range("F32").Value = 469
$ws4.Range("F34").Value = 73
$ws4.Range("F35").Value = 1610
$ws4.Range("F36").Value = 244
$ws4.Range("F37").Value = 443
$ws4.Range("G37").Value = "不可售"
$ws4.Range("F39").Value = 133
$ws4.Range("F40").Value = 453
$ws4.Range("F41").Value = 7
$ws4.Range("F42").Value = 109
$ws4.Range("F44").Value = 121
$ws4.Range("F49").Value = 20
